# Insert a new data row above row 490, shifting existing rows 490:584 down to 491:585.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("490:490").Insert()

# Populate the newly inserted row 490 with its values.
$ws.Range("A490").Value = 4
$ws.Range("B490").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C490").Value = "Los Lagos"
$ws.Range("D490").Value = 45275
$ws.Range("E490").Value = 10
$ws.Range("F490").Value = 100112045
$ws.Range("G490").Value = "Zapallo"
$ws.Range("H490").Value = "Paine"
$ws.Range("I490").Value = "1a nueva(o)"
$ws.Range("J490").Value = 1500
$ws.Range("K490").Value = 1100
$ws.Range("L490").Value = 1200
$ws.Range("M490").Value = 1150
$ws.Range("N490").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O490").Value = "Región de O'Higgins"
$ws.Range("P490").Value = 1150
$ws.Range("Q490").Value = 1
$ws.Range("R490").Value = "Hortaliza"

# Make sure the date cell keeps the expected date format used throughout column D.
$ws.Range("D490").NumberFormat = "YYYY-MM-DD HH:MM:SS"
